# Work diary for the morning
# Adds a new entry to the "JdT-TPI_LRD" work-log table: 12/05/2022,
# Réalisation, 3h, "Vérification de l'email de l'utilisateur".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date column (A44) was stored with a custom "dd/mm/yyyy" format;
# normalise it to the plain short-date builtin format before extending
# the table, matching the look of the rest of the date column.
$ws.Range("A44").NumberFormat = "mm-dd-yy"

# Grow the "Tableau1" table by one row so the new entry inherits the
# table's styling/autofilter/formula range automatically.
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

$newRow.Range.Cells(1, 1).Value = 44693
$ws.Range("A44").Copy()
$newRow.Range.Cells(1, 1).PasteSpecial(-4122)
$newRow.Range.Cells(1, 2).Value = "Réalisation"
$newRow.Range.Cells(1, 3).Value = 3
$newRow.Range.Cells(1, 4).Value = "Vérification de l'email de l'utilisateur"

# Match the author's final selection after typing the new row.
$ws.Range("F45").Select()
